$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 16.717
$ws.Range("C9").Value = -11.292
$ws.Range("E12").Value = 17.503
$ws.Range("C13").Value = -13.361
$ws.Range("E14").Value = 16.68
$ws.Range("C16").Value = -12.476
$ws.Range("C18").Value = -11.62
$ws.Range("E19").Value = 16.472
$ws.Range("C20").Value = -12.417
$ws.Range("C26").Value = -12.405
$ws.Range("E26").Value = 16.278
$ws.Range("C27").Value = -13.533
$ws.Range("E27").Value = 16.509
$ws.Range("C29").Value = -12.434
$ws.Range("E29").Value = 16.856
$ws.Range("C35").Value = -12.377
$ws.Range("C36").Value = -12.428
$ws.Range("E37").Value = 16.813
$ws.Range("E38").Value = 16.741
$ws.Range("C45").Value = -13.006
$ws.Range("E47").Value = 16.69
$ws.Range("E51").Value = 16.867
$ws.Range("E52").Value = 16.955
$ws.Range("C55").Value = -13.381
$ws.Range("E55").Value = 16.25
$ws.Range("C57").Value = -13.426
$ws.Range("C69").Value = -10.672
$ws.Range("E69").Value = 17.396
$ws.Range("E70").Value = 17.503
$ws.Range("C76").Value = -13.094
$ws.Range("E76").Value = 16.764
$ws.Range("C78").Value = -12.5
$ws.Range("E81").Value = 16.206
$ws.Range("C82").Value = -12.323
$ws.Range("C83").Value = -13.064
$ws.Range("E83").Value = 16.731
$ws.Range("C93").Value = -11.587
$ws.Range("E94").Value = 18.06
$ws.Range("C97").Value = -12.208
$ws.Range("E100").Value = 16.634
$ws.Range("E102").Value = 16.745
